$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 216.23077
$ws.Range("I33").Value = 116.42857
$ws.Range("J33").Value = 332.66666
$ws.Range("K33").Value = 116.42857
$ws.Range("L33").Value = 332.66666
$ws.Range("M33").Value = 112.57143
$ws.Range("N33").Value = -790.66666

$ws.Range("H92").Value = 551.2
$ws.Range("I92").Value = 551.2
$ws.Range("K92").Value = 551.2
$ws.Range("M92").Value = 696.8

$ws.Range("H98").Value = 6389.276
$ws.Range("I98").Value = 6590.5386
$ws.Range("J98").Value = 4645
$ws.Range("K98").Value = 6590.5386
$ws.Range("L98").Value = 4645
$ws.Range("M98").Value = -5092.5386
$ws.Range("N98").Value = -7641

$ws.Range("H100").Value = 1923.8182
$ws.Range("I100").Value = 1991.2
$ws.Range("K100").Value = 1991.2
$ws.Range("M100").Value = -1450.2

$ws.Range("H116").Value = 4193.8
$ws.Range("I116").Value = 3792.25
$ws.Range("K116").Value = 3792.25
$ws.Range("M116").Value = -350.25

$ws.Range("H122").Value = 6389.276
$ws.Range("I122").Value = 6590.5386
$ws.Range("J122").Value = 4645
$ws.Range("K122").Value = 19771.6158
$ws.Range("L122").Value = 13935
$ws.Range("M122").Value = -17321.6158
$ws.Range("N122").Value = -18835

$ws.Range("H132").Value = 8053.8667
$ws.Range("I132").Value = 8414.929
$ws.Range("K132").Value = 25244.787
$ws.Range("M132").Value = -22714.787

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 52030.24
$ws.Range("I45").Value = 80165.96000000001
$ws.Range("K45").Value = 80165.96000000001
$ws.Range("M45").Value = -79788.96000000001

$ws.Range("H61").Value = 3274246.8
$ws.Range("I61").Value = 4391022
$ws.Range("K61").Value = 4391022
$ws.Range("M61").Value = -4390810

$ws.Range("H63").Value = 2002804.4
$ws.Range("I63").Value = 2995.6667
$ws.Range("K63").Value = 2995.6667
$ws.Range("M63").Value = -2309.6667

$ws.Range("H66").Value = 2002804.4
$ws.Range("I66").Value = 2995.6667
$ws.Range("K66").Value = 14978.3335
$ws.Range("M66").Value = -11546.3335

$ws.Range("H74").Value = 3725.394
$ws.Range("I74").Value = 2289.125
$ws.Range("K74").Value = 2289.125
$ws.Range("M74").Value = -1415.125

$ws.Range("H77").Value = 3725.394
$ws.Range("I77").Value = 2289.125
$ws.Range("K77").Value = 11445.625
$ws.Range("M77").Value = -7077.625

$ws.Range("H97").Value = 530.2941
$ws.Range("I97").Value = 482.1875
$ws.Range("K97").Value = 482.1875
$ws.Range("M97").Value = 13.8125

$ws.Range("H132").Value = 4032.4243
$ws.Range("I132").Value = 3387.157
$ws.Range("K132").Value = 10161.471
$ws.Range("M132").Value = -7631.471000000001

$ws.Range("H136").Value = 3274246.8
$ws.Range("I136").Value = 4391022
$ws.Range("K136").Value = 13173066
$ws.Range("M136").Value = -13170516

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 100000
$ws.Range("J132").Value = 100000
$ws.Range("L132").Value = 100000
$ws.Range("N132").Value = -110120

$ws.Range("H134").Value = 5100.3784
$ws.Range("I134").Value = 5167.091
$ws.Range("K134").Value = 15501.273
$ws.Range("M134").Value = -12966.273

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5233.5806
$ws.Range("I31").Value = 4650.875
$ws.Range("J31").Value = 5855.1333
$ws.Range("K31").Value = 4650.875
$ws.Range("L31").Value = 5855.1333
$ws.Range("M31").Value = -4355.875
$ws.Range("N31").Value = -6445.1333

$ws.Range("H34").Value = 5233.5806
$ws.Range("I34").Value = 4650.875
$ws.Range("J34").Value = 5855.1333
$ws.Range("K34").Value = 4650.875
$ws.Range("L34").Value = 5855.1333
$ws.Range("M34").Value = -4448.875
$ws.Range("N34").Value = -6259.1333

$ws.Range("H105").Value = 1294.7142
$ws.Range("I105").Value = 1299.421
$ws.Range("K105").Value = 1299.421
$ws.Range("M105").Value = 447.579

$ws.Range("H132").Value = 2017
$ws.Range("I132").Value = 2017
$ws.Range("K132").Value = 6051
$ws.Range("M132").Value = -3521

$ws.Range("H139").Value = 139912.5
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 139912.5
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 139912.5
$ws.Range("M139").ClearContents()
$ws.Range("N139").Value = -150192.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2919.1667
$ws.Range("I5").Value = 1980.25
$ws.Range("J5").Value = 4797
$ws.Range("K5").Value = 5940.75
$ws.Range("L5").Value = 14391
$ws.Range("M5").Value = -5828.75
$ws.Range("N5").Value = -14615

$ws.Range("H11").Value = 8400290
$ws.Range("I11").Value = 8400290
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 25200870
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -25200730
$ws.Range("N11").ClearContents()

$ws.Range("H75").Value = 900
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()

$ws.Range("H78").Value = 900
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()

$ws.Range("H135").Value = 2919.1667
$ws.Range("I135").Value = 1980.25
$ws.Range("J135").Value = 4797
$ws.Range("K135").Value = 17822.25
$ws.Range("L135").Value = 43173
$ws.Range("M135").Value = -15287.25
$ws.Range("N135").Value = -48243

$ws.Range("H140").Value = 3168.45
$ws.Range("J140").Value = 3075
$ws.Range("L140").Value = 9225
$ws.Range("N140").Value = -19585

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4703.421
$ws.Range("I136").Value = 4843.375
$ws.Range("K136").Value = 14530.125
$ws.Range("M136").Value = -11980.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()

$ws.Range("H113").Value = 376.61905
$ws.Range("I113").Value = 385.21054
$ws.Range("K113").Value = 1155.63162
$ws.Range("M113").Value = 1014.36838

$ws.Range("H117").Value = 78000
$ws.Range("J117").Value = 78000
$ws.Range("L117").Value = 78000
$ws.Range("N117").Value = -87178

$ws.Range("H122").Value = 3585.0322
$ws.Range("I122").Value = 2416.818
$ws.Range("K122").Value = 7250.454000000001
$ws.Range("M122").Value = -4800.454000000001

$ws.Range("H132").Value = 5034.217
$ws.Range("I132").Value = 4531.735
$ws.Range("K132").Value = 13595.205
$ws.Range("M132").Value = -11065.205
